$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Ph.D. thesis #15, "Overfitting in Automated Program Repair: Challenges
# and Solutions") previously cited the "Patch correctness assessment bias 3/4"
# as the bias it covered, and had no entry for "omitted bias". The commit
# corrects this to the actual bias reported by the thesis ("Only-manual
# validation bias" / "Only-independent test validation bias") and marks the
# omitted-bias column as "No".
$ws.Range("D16").Value = 'It reported the "Only-manual validation bias" and "Only-independent test validation bias", the same biases reported by Le et al. \cite{le2019reliability}.'
$ws.Range("E16").Value = "No"

# Reflect the author's final cursor/view position: scrolled down so row 7 is
# the first visible row, with D12 selected.
$win = $excel.ActiveWindow
try { $win.ScrollRow = 7 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("D12").Select()
